$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(23).Insert()

$ws.Cells.Item(23, 1).Value = "x-ext"
$ws.Cells.Item(23, 2).Value = "Race Code (Adams County)"
$ws.Cells.Item(23, 3).Value = "Person Race"
$ws.Cells.Item(23, 5).Value = "/br-doc:BookingReport/nc:Person[@structures:id=/br-doc:BookingReport/j:Booking/j:BookingSubject/nc:RoleOfPerson/@structures:ref]/ac-bkg-codes:PersonRaceCode"
